# Scheduled-runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# figures (columns H-N) across the ALC/ARM/BSM/CRP/CUL/LTW/WVR leve tables
# following an upstream market-price data pull. GSM is unaffected this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2857.4075  # ALC!H40
$ws.Cells.Item(40, 9).Value = 3563.4546  # ALC!I40
$ws.Cells.Item(40, 10).Value = 2372  # ALC!J40
$ws.Cells.Item(40, 11).Value = 3563.4546  # ALC!K40
$ws.Cells.Item(40, 12).Value = 2372  # ALC!L40
$ws.Cells.Item(40, 13).Value = -3388.4546  # ALC!M40
$ws.Cells.Item(40, 14).Value = -2722  # ALC!N40

$ws.Cells.Item(64, 8).Value = 2825  # ALC!H64
$ws.Cells.Item(64, 9).Value = 2677.889  # ALC!I64
$ws.Cells.Item(64, 10).Value = 3089.8  # ALC!J64
$ws.Cells.Item(64, 11).Value = 2677.889  # ALC!K64
$ws.Cells.Item(64, 12).Value = 3089.8  # ALC!L64
$ws.Cells.Item(64, 13).Value = -2429.889  # ALC!M64
$ws.Cells.Item(64, 14).Value = -3585.8  # ALC!N64

$ws.Cells.Item(67, 8).Value = 2825  # ALC!H67
$ws.Cells.Item(67, 9).Value = 2677.889  # ALC!I67
$ws.Cells.Item(67, 10).Value = 3089.8  # ALC!J67
$ws.Cells.Item(67, 11).Value = 2677.889  # ALC!K67
$ws.Cells.Item(67, 12).Value = 3089.8  # ALC!L67
$ws.Cells.Item(67, 13).Value = -1819.889  # ALC!M67
$ws.Cells.Item(67, 14).Value = -4805.8  # ALC!N67

$ws.Cells.Item(93, 8).Value = 35267.668  # ALC!H93
$ws.Cells.Item(93, 10).Value = 35267.668  # ALC!J93
$ws.Cells.Item(93, 12).Value = 35267.668  # ALC!L93
$ws.Cells.Item(93, 14).Value = -40259.668  # ALC!N93

$ws.Cells.Item(138, 8).Value = 5210563.5  # ALC!H138
$ws.Cells.Item(138, 9).Value = 1172.6216  # ALC!I138
$ws.Cells.Item(138, 10).Value = 22733060  # ALC!J138
$ws.Cells.Item(138, 11).Value = 3517.8648  # ALC!K138
$ws.Cells.Item(138, 12).Value = 68199180  # ALC!L138
$ws.Cells.Item(138, 13).Value = 1622.1352  # ALC!M138
$ws.Cells.Item(138, 14).Value = -68209460  # ALC!N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(94, 8).Value = 50330  # ARM!H94
$ws.Cells.Item(94, 10).Value = 50330  # ARM!J94
$ws.Cells.Item(94, 12).Value = 50330  # ARM!L94
$ws.Cells.Item(94, 14).Value = -52132  # ARM!N94

$ws.Cells.Item(103, 8).Value = 100000000  # ARM!H103
$ws.Cells.Item(103, 10).Value = 100000000  # ARM!J103
$ws.Cells.Item(103, 12).Value = 100000000  # ARM!L103
$ws.Cells.Item(103, 14).Value = -100002344  # ARM!N103

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3414.2407  # BSM!H134
$ws.Cells.Item(134, 9).Value = 2502.2092  # BSM!I134
$ws.Cells.Item(134, 10).Value = 6979.4546  # BSM!J134
$ws.Cells.Item(134, 11).Value = 7506.6276  # BSM!K134
$ws.Cells.Item(134, 12).Value = 20938.3638  # BSM!L134
$ws.Cells.Item(134, 13).Value = -4971.6276  # BSM!M134
$ws.Cells.Item(134, 14).Value = -26008.3638  # BSM!N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(36, 8).Value = 1952.4  # CRP!H36
$ws.Cells.Item(36, 9).Value = 1952.4  # CRP!I36
$ws.Cells.Item(36, 11).Value = 1952.4  # CRP!K36
$ws.Cells.Item(36, 13).Value = -1564.4  # CRP!M36

$ws.Cells.Item(40, 8).Value = 1952.4  # CRP!H40
$ws.Cells.Item(40, 9).Value = 1952.4  # CRP!I40
$ws.Cells.Item(40, 11).Value = 1952.4  # CRP!K40
$ws.Cells.Item(40, 13).Value = -1792.4  # CRP!M40

$ws.Cells.Item(70, 8).Value = 30516.666  # CRP!H70
$ws.Cells.Item(70, 10).Value = 30516.666  # CRP!J70
$ws.Cells.Item(70, 12).Value = 30516.666  # CRP!L70
$ws.Cells.Item(70, 14).Value = -31146.666  # CRP!N70

$ws.Cells.Item(73, 8).Value = 30516.666  # CRP!H73
$ws.Cells.Item(73, 10).Value = 30516.666  # CRP!J73
$ws.Cells.Item(73, 12).Value = 30516.666  # CRP!L73
$ws.Cells.Item(73, 14).Value = -32700.666  # CRP!N73

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 261  # CUL!H26
$ws.Cells.Item(26, 9).Value = 168.625  # CUL!I26
$ws.Cells.Item(26, 10).Value = 1000  # CUL!J26
$ws.Cells.Item(26, 11).Value = 505.875  # CUL!K26
$ws.Cells.Item(26, 12).Value = 3000  # CUL!L26
$ws.Cells.Item(26, 13).Value = -217.875  # CUL!M26
$ws.Cells.Item(26, 14).Value = -3576  # CUL!N26

$ws.Cells.Item(62, 8).Value = 5320.8823  # CUL!H62
$ws.Cells.Item(62, 9).Value = 999.6667  # CUL!I62
$ws.Cells.Item(62, 10).Value = 6246.857  # CUL!J62
$ws.Cells.Item(62, 11).Value = 2999.0001  # CUL!K62
$ws.Cells.Item(62, 12).Value = 18740.571  # CUL!L62
$ws.Cells.Item(62, 13).Value = -2313.0001  # CUL!M62
$ws.Cells.Item(62, 14).Value = -20112.571  # CUL!N62

$ws.Cells.Item(63, 8).Value = 2971  # CUL!H63
$ws.Cells.Item(63, 9).Value = 618.3333  # CUL!I63
$ws.Cells.Item(63, 10).Value = 6500  # CUL!J63
$ws.Cells.Item(63, 11).Value = 1854.9999  # CUL!K63
$ws.Cells.Item(63, 12).Value = 19500  # CUL!L63
$ws.Cells.Item(63, 13).Value = -1105.9999  # CUL!M63
$ws.Cells.Item(63, 14).Value = -20998  # CUL!N63

$ws.Cells.Item(64, 8).Value = 3142.8572  # CUL!H64

$ws.Cells.Item(65, 8).Value = 5320.8823  # CUL!H65
$ws.Cells.Item(65, 9).Value = 999.6667  # CUL!I65
$ws.Cells.Item(65, 10).Value = 6246.857  # CUL!J65
$ws.Cells.Item(65, 11).Value = 8997.0003  # CUL!K65
$ws.Cells.Item(65, 12).Value = 56221.713  # CUL!L65
$ws.Cells.Item(65, 13).Value = -5565.0003  # CUL!M65
$ws.Cells.Item(65, 14).Value = -63085.713  # CUL!N65

$ws.Cells.Item(66, 8).Value = 2971  # CUL!H66
$ws.Cells.Item(66, 9).Value = 618.3333  # CUL!I66
$ws.Cells.Item(66, 10).Value = 6500  # CUL!J66
$ws.Cells.Item(66, 11).Value = 5564.9997  # CUL!K66
$ws.Cells.Item(66, 12).Value = 58500  # CUL!L66
$ws.Cells.Item(66, 13).Value = -1820.9997  # CUL!M66
$ws.Cells.Item(66, 14).Value = -65988  # CUL!N66

$ws.Cells.Item(67, 8).Value = 3142.8572  # CUL!H67

$ws.Cells.Item(74, 8).Value = 8028.125  # CUL!H74
$ws.Cells.Item(74, 9).Value = 3306.5  # CUL!I74
$ws.Cells.Item(74, 10).Value = 12749.75  # CUL!J74
$ws.Cells.Item(74, 11).Value = 9919.5  # CUL!K74
$ws.Cells.Item(74, 12).Value = 38249.25  # CUL!L74
$ws.Cells.Item(74, 13).Value = -8858.5  # CUL!M74
$ws.Cells.Item(74, 14).Value = -40371.25  # CUL!N74

$ws.Cells.Item(77, 8).Value = 8028.125  # CUL!H77
$ws.Cells.Item(77, 9).Value = 3306.5  # CUL!I77
$ws.Cells.Item(77, 10).Value = 12749.75  # CUL!J77
$ws.Cells.Item(77, 11).Value = 29758.5  # CUL!K77
$ws.Cells.Item(77, 12).Value = 114747.75  # CUL!L77
$ws.Cells.Item(77, 13).Value = -24454.5  # CUL!M77
$ws.Cells.Item(77, 14).Value = -125355.75  # CUL!N77

$ws.Cells.Item(132, 8).Value = 6409.091  # CUL!H132
$ws.Cells.Item(132, 9).Value = 0  # CUL!I132
$ws.Cells.Item(132, 10).Value = 6409.091  # CUL!J132
$ws.Cells.Item(132, 11).Value = 0  # CUL!K132
$ws.Cells.Item(132, 12).Value = 57681.819  # CUL!L132
$ws.Cells.Item(132, 13).ClearContents()  # CUL!M132: remove (was -4744.454299999999)
$ws.Cells.Item(132, 14).Value = -62741.819  # CUL!N132

$ws.Cells.Item(141, 8).Value = 1532.5  # CUL!H141
$ws.Cells.Item(141, 9).Value = 1532.5  # CUL!I141
$ws.Cells.Item(141, 11).Value = 4597.5  # CUL!K141
$ws.Cells.Item(141, 13).Value = 582.5  # CUL!M141

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 4059.6428  # LTW!H16
$ws.Cells.Item(16, 9).Value = 2893.5  # LTW!I16
$ws.Cells.Item(16, 10).Value = 6975  # LTW!J16
$ws.Cells.Item(16, 11).Value = 2893.5  # LTW!K16
$ws.Cells.Item(16, 12).Value = 6975  # LTW!L16
$ws.Cells.Item(16, 13).Value = -2723.5  # LTW!M16
$ws.Cells.Item(16, 14).Value = -7315  # LTW!N16

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 10239.667  # WVR!H41
$ws.Cells.Item(41, 9).Value = 8342  # WVR!I41
$ws.Cells.Item(41, 11).Value = 8342  # WVR!K41
$ws.Cells.Item(41, 13).Value = -7952  # WVR!M41

$ws.Cells.Item(48, 8).Value = 0  # WVR!H48
$ws.Cells.Item(48, 9).Value = 0  # WVR!I48
$ws.Cells.Item(48, 10).Value = 0  # WVR!J48
$ws.Cells.Item(48, 11).Value = 0  # WVR!K48
$ws.Cells.Item(48, 12).Value = 0  # WVR!L48
$ws.Cells.Item(48, 13).ClearContents()  # WVR!M48: remove (was -10931)
$ws.Cells.Item(48, 14).ClearContents()  # WVR!N48: remove (was -21138)

$ws.Cells.Item(95, 8).Value = 50344  # WVR!H95
$ws.Cells.Item(95, 10).Value = 50344  # WVR!J95
$ws.Cells.Item(95, 12).Value = 50344  # WVR!L95
$ws.Cells.Item(95, 14).Value = -55836  # WVR!N95

$ws.Cells.Item(107, 8).Value = 1984.375  # WVR!H107
$ws.Cells.Item(107, 9).Value = 2362  # WVR!I107
$ws.Cells.Item(107, 10).Value = 851.5  # WVR!J107
$ws.Cells.Item(107, 11).Value = 7086  # WVR!K107
$ws.Cells.Item(107, 12).Value = 2554.5  # WVR!L107
$ws.Cells.Item(107, 13).Value = -5166  # WVR!M107
$ws.Cells.Item(107, 14).Value = -6394.5  # WVR!N107

$ws.Cells.Item(136, 8).Value = 3572514  # WVR!H136
$ws.Cells.Item(136, 9).Value = 4546509  # WVR!I136
$ws.Cells.Item(136, 10).Value = 1198.3334  # WVR!J136
$ws.Cells.Item(136, 11).Value = 13639527  # WVR!K136
$ws.Cells.Item(136, 12).Value = 3595.0002  # WVR!L136
$ws.Cells.Item(136, 13).Value = -13636977  # WVR!M136
$ws.Cells.Item(136, 14).Value = -8695.0002  # WVR!N136
